# Populate Sheet1 with the Issue Report demo data
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header row
$ws.Range("A1").Value = "number"
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "description"
$ws.Range("D1").Value = "project"
$ws.Range("E1").Value = "type"
$ws.Range("F1").Value = "corrective_action"
$ws.Range("G1").Value = "cost"
$ws.Range("H1").Value = "status"
$ws.Range("I1").Value = "create_time"

# Row 2 - IR00001
$ws.Range("A2").Value = "IR00001"
$ws.Range("B2").Value = "Incorrect PCB layout"
$ws.Range("C2").Value = "Connectivity of PCB-type inductor in layout is incorrect."
$ws.Range("D2").Value = "P09-001"
$ws.Range("E2").Value = "PCB"
$ws.Range("F2").Value = "YES"
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = "CLOSED"
$ws.Range("I2").Value = 38362

# Row 3 - IR00002
$ws.Range("A3").Value = "IR00002"
$ws.Range("B3").Value = "Difficult assembly procedure"
$ws.Range("C3").Value = "Difficult to assemble and test PCA due to copper-side-up orientation. Orienting PCA component-side up allows for full access to components for debugging after assembly."
$ws.Range("D3").Value = "P09-001"
$ws.Range("E3").Value = "ASSY"
$ws.Range("F3").Value = "YES"
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = "CLOSED"
$ws.Range("I3").Value = 38407

# Wrap text on the description column (applied first so it becomes style index 1)
$ws.Range("C1:C3").WrapText = $true

# Date formatting for the create_time column (built-in date format -> style index 2)
$ws.Range("I2").NumberFormat = "mm-dd-yy"
$ws.Range("I3").NumberFormat = "mm-dd-yy"

# Column widths (bestFit results from the original workbook, reproduced as closely
# as this engine's char-width quantization allows)
$ws.Range("A:A").ColumnWidth = 7.16666666666667
$ws.Range("B:B").ColumnWidth = 26.1666666666667
$ws.Range("C:C").ColumnWidth = 38.7369791666667
$ws.Range("D:D").ColumnWidth = 6.45182291666667
$ws.Range("E:E").ColumnWidth = 4.16666666666667
$ws.Range("F:F").ColumnWidth = 15.5924479166667
$ws.Range("G:G").ColumnWidth = 3.73697916666667
$ws.Range("H:H").ColumnWidth = 6.87760416666667
$ws.Range("I:I").ColumnWidth = 10.8776041666667

# Row heights to match the wrapped text rows
$ws.Range("2:2").RowHeight = 30
$ws.Range("3:3").RowHeight = 75
